# Weekly data refresh: insert a new row for the latest week (Red Globe,
# Provincia del Elquí) at the top of the data block and shift the
# existing rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(86).Insert()

$ws.Range("A86").Value = 4
$ws.Range("B86").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C86").Value = "Los Lagos"
$ws.Range("D86").Value = 44586
$ws.Range("D86").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E86").Value = 10
$ws.Range("F86").Value = "Fruta"
$ws.Range("G86").Value = 100109
$ws.Range("H86").Value = "Uva"
$ws.Range("I86").Value = 100109001
$ws.Range("J86").Value = "Uva"
$ws.Range("K86").Value = "Red Globe"
$ws.Range("L86").Value = "Primera"
$ws.Range("M86").Value = 300
$ws.Range("N86").Value = 14000
$ws.Range("O86").Value = 14000
$ws.Range("P86").Value = 14000
$ws.Range("Q86").Value = "`$/bandeja 12 kilos"
$ws.Range("R86").Value = "Provincia del Elquí"
$ws.Range("S86").Value = 1167
$ws.Range("T86").Value = 12
